$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "mercredi / FSQTEL - C" block that used to sit under
# "Semaine : 48" (rows 22-24). This shifts the following "vendredi / FSQTEL"
# (controle) block up to become rows 22-24, matching the new week 48 content,
# and shrinks the used range from A1:I27 to A1:I24.
$ws.Rows.Item(22).Resize(3).Delete()

# Fill in the room ("salle") column (F) for every course/TD/controle line.
$ws.Range("F3").Value = "U3-Amphi"
$ws.Range("F4").Value = "U3-Amphi"
$ws.Range("F7").Value = "U3-Amphi"
$ws.Range("F8").Value = "U3-Amphi"
$ws.Range("F11").Value = "U3-Amphi"
$ws.Range("F12").Value = "U3-Amphi"
$ws.Range("F15").Value = "U3-110"
$ws.Range("F16").Value = "U3-110"
$ws.Range("F19").Value = "U3-4"
$ws.Range("F20").Value = "U3-Amphi"
$ws.Range("F23").Value = "U3-Amphi"
$ws.Range("F24").Value = "U3-Amphi"
